# This script updates the "AI_statistics" summary-statistics sheet
# (count/mean/std/25%/50%/75% rows for the 18 survey items) to reflect the
# refreshed data snapshot referenced by the commit "Kinda final results I guess?"
# Only cell values change; no structural or style changes are required.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 34
$ws.Range("C2").Value = 34
$ws.Range("D2").Value = 34
$ws.Range("E2").Value = 34
$ws.Range("F2").Value = 34
$ws.Range("G2").Value = 34
$ws.Range("H2").Value = 34
$ws.Range("I2").Value = 34
$ws.Range("J2").Value = 34
$ws.Range("K2").Value = 34
$ws.Range("L2").Value = 34
$ws.Range("M2").Value = 34
$ws.Range("N2").Value = 34
$ws.Range("O2").Value = 34
$ws.Range("P2").Value = 34
$ws.Range("Q2").Value = 34
$ws.Range("R2").Value = 34
$ws.Range("S2").Value = 34
$ws.Range("T2").Value = 34
$ws.Range("B3").Value = 3.264705882352941
$ws.Range("C3").Value = 3.5
$ws.Range("D3").Value = 3.029411764705882
$ws.Range("E3").Value = 2.794117647058823
$ws.Range("F3").Value = 3.058823529411764
$ws.Range("G3").Value = 3.352941176470588
$ws.Range("H3").Value = 2.705882352941177
$ws.Range("I3").Value = 2.852941176470588
$ws.Range("J3").Value = 2.823529411764706
$ws.Range("K3").Value = 3.352941176470588
$ws.Range("L3").Value = 3.323529411764706
$ws.Range("M3").Value = 2.823529411764706
$ws.Range("N3").Value = 3.147058823529412
$ws.Range("O3").Value = 3.676470588235294
$ws.Range("P3").Value = 3.205882352941177
$ws.Range("Q3").Value = 3.382352941176471
$ws.Range("R3").Value = 3.323529411764706
$ws.Range("S3").Value = 3.411764705882353
$ws.Range("T3").Value = 2.735294117647059
$ws.Range("B4").Value = 1.377499348939224
$ws.Range("C4").Value = 1.308480497417219
$ws.Range("D4").Value = 1.466500557762595
$ws.Range("E4").Value = 1.665685985989427
$ws.Range("F4").Value = 1.347077115755897
$ws.Range("G4").Value = 1.453988306842575
$ws.Range("H4").Value = 1.54781197990379
$ws.Range("I4").Value = 1.479810654169145
$ws.Range("J4").Value = 1.381053358091791
$ws.Range("K4").Value = 1.432995441663677
$ws.Range("L4").Value = 1.248528545693596
$ws.Range("M4").Value = 1.381053358091791
$ws.Range("N4").Value = 1.43827238414022
$ws.Range("O4").Value = 1.173458711499294
$ws.Range("P4").Value = 1.365803388057981
$ws.Range("Q4").Value = 1.181029491391532
$ws.Range("R4").Value = 1.429570569877629
$ws.Range("S4").Value = 1.328422328310143
$ws.Range("T4").Value = 1.54348726628258
$ws.Range("I6").Value = 2
$ws.Range("Q6").Value = 2.25
$ws.Range("B7").Value = 4
$ws.Range("E7").Value = 2.5
$ws.Range("G7").Value = 4
$ws.Range("I7").Value = 2.5
$ws.Range("L7").Value = 3.5
$ws.Range("M7").Value = 2.5
$ws.Range("R7").Value = 3.5
$ws.Range("C8").Value = 4.75
$ws.Range("E8").Value = 4.75
$ws.Range("N8").Value = 4.75
$ws.Range("O8").Value = 4.75
$ws.Range("P8").Value = 4
$ws.Range("S8").Value = 4.75
